$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text-looking numbers stored as inline strings in the
# original workbook (t="inlineStr"). Excel's COM layer auto-converts plain
# numeric-looking text into real numbers, so we force a text number format
# before writing the value and then restore the "Normal" style so no stray
# style index is left behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "250.30"
Set-TextValue $ws.Range("D3") "24.06"
Set-TextValue $ws.Range("D4") "6.034"
Set-TextValue $ws.Range("D5") "0.05987"
Set-TextValue $ws.Range("D6") "3.429"
Set-TextValue $ws.Range("D7") "6.565"
Set-TextValue $ws.Range("D8") "1.320"
Set-TextValue $ws.Range("D9") "0.7996"
Set-TextValue $ws.Range("D10") "0.1491"
Set-TextValue $ws.Range("D11") "0.07930"
Set-TextValue $ws.Range("D12") "0.03349"
Set-TextValue $ws.Range("D14") "0.09276"
Set-TextValue $ws.Range("D15") "3.568"
Set-TextValue $ws.Range("D16") "0.001682"
Set-TextValue $ws.Range("D17") "0.04779"
Set-TextValue $ws.Range("D18") "0.0006100"
Set-TextValue $ws.Range("D19") "0.006235"
Set-TextValue $ws.Range("D20") "0.005695"
Set-TextValue $ws.Range("D21") "0.001073"
Set-TextValue $ws.Range("D22") "0.0001504"
Set-TextValue $ws.Range("D23") "3.676"
Set-TextValue $ws.Range("D24") "2.203"
Set-TextValue $ws.Range("D27") "0.0006496"
Set-TextValue $ws.Range("D40") "0.04454"
Set-TextValue $ws.Range("D41") "0.007062"

# Rows 42 and 43 swap their coin identity (BKEXToken <-> CEJI)
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D42") "0.003610"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D43") "0.1070"
$ws.Range("E43").Value = "42BKEXTokenBKK"

Set-TextValue $ws.Range("D44") "0.01027"
Set-TextValue $ws.Range("D45") "0.002468"
Set-TextValue $ws.Range("D46") "0.00005901"
Set-TextValue $ws.Range("D47") "0.00000000752"
Set-TextValue $ws.Range("D48") "0.7024"
Set-TextValue $ws.Range("D49") "0.1176"
Set-TextValue $ws.Range("D50") "0.00002106"
Set-TextValue $ws.Range("D51") "0.01013"
